# "Added All commodity and validation"
# - Extend "wheat" and "rra" sheets from 3 columns (From/To/Commodity) to
#   6 columns (From/From State/To/To State/Commodity/Values).
# - Remove the sample data row that used to live on "rra".
# - Add six new sheets (coarse_grain, frk_rra, frk_br, frk, frkcgr, wcgr),
#   each with the same 6-column header; frk_rra additionally ships three
#   data rows.

$wb = $excel.ActiveWorkbook

function Set-DptHeader($ws) {
    # Make room for the 3 new columns by copying the existing formatted
    # header cell across first, so the new cells inherit the bold /
    # bordered / centered header look -- then stamp the real text in.
    $ws.Range("C1").Copy($ws.Range("D1:F1"))
    $ws.Range("A1").Value = "From"
    $ws.Range("B1").Value = "From State"
    $ws.Range("C1").Value = "To"
    $ws.Range("D1").Value = "To State"
    $ws.Range("E1").Value = "Commodity"
    $ws.Range("F1").Value = "Values"
}

# ---------------------------------------------------------------------
# 1) "wheat": From/To/Commodity -> From/From State/To/To State/Commodity/Values
# ---------------------------------------------------------------------
$wheat = $wb.Worksheets.Item("wheat")
Set-DptHeader($wheat)

# ---------------------------------------------------------------------
# 2) "rra": same header upgrade, and drop the old sample data row so only
#    the header row remains.
# ---------------------------------------------------------------------
$rra = $wb.Worksheets.Item("rra")
$rra.Rows.Item(2).Delete()
Set-DptHeader($rra)

# ---------------------------------------------------------------------
# 3) Six brand-new sheets, appended after "rra", each starting out with
#    the same 6-column header (copied wholesale from the now-fixed
#    "wheat" sheet so formatting matches exactly).
# ---------------------------------------------------------------------
$newSheetNames = @("coarse_grain", "frk_rra", "frk_br", "frk", "frkcgr", "wcgr")
$createdSheets = @{}
foreach ($name in $newSheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
    $wheat.Range("A1:F1").Copy($newSheet.Range("A1:F1"))
    $createdSheets[$name] = $newSheet
}

# ---------------------------------------------------------------------
# 4) "frk_rra" additionally gets three data rows.
# ---------------------------------------------------------------------
$frkRra = $createdSheets["frk_rra"]

$frkRra.Range("A2").Value = "BIDR"
$frkRra.Range("B2").Value = "Karnataka"
$frkRra.Range("C2").Value = "BBMN"
$frkRra.Range("D2").Value = "Jammu & Kashmir"
$frkRra.Range("E2").Value = "FRK RRA"
$frkRra.Range("F2").Value = 1

$frkRra.Range("A3").Value = "BUDI"
$frkRra.Range("B3").Value = "Rajasthan"
$frkRra.Range("C3").Value = "BRW"
$frkRra.Range("D3").Value = "Punjab"
$frkRra.Range("E3").Value = "FRK RRA"
$frkRra.Range("F3").Value = 1

$frkRra.Range("A4").Value = "CHD"
$frkRra.Range("B4").Value = "MP"
$frkRra.Range("C4").Value = "HZBN"
$frkRra.Range("D4").Value = "Jharkhand"
$frkRra.Range("E4").Value = "FRK RRA"
$frkRra.Range("F4").Value = 1
